$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing daily sales values (column B) for a few days in August (rows 9, 11, 13, 14)
$ws.Range("B9").Value = 117963.69
$ws.Range("B11").Value = 14233.4
$ws.Range("B13").Value = 11107.75
$ws.Range("B14").Value = 11839.02

# Insert two new rows at position 15 to hold new daily records (days 20 and 21 of August/2025),
# pushing the existing rows (old row 15 onward) down by two rows.
$ws.Rows.Item(15).Resize(2).Insert()

# Fill in the two newly inserted rows with the new data
$ws.Range("A15").Value = 20
$ws.Range("B15").Value = 22359.16
$ws.Range("C15").Value = 8
$ws.Range("D15").Value = 2025
$ws.Range("E15").Value = "08/2025"

$ws.Range("A16").Value = 21
$ws.Range("B16").Value = 18859.31
$ws.Range("C16").Value = 8
$ws.Range("D16").Value = 2025
$ws.Range("E16").Value = "08/2025"
